$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-08-09 Saturday" "2025-08-10 Sunday"

Replace-Text "39÷9=" "33÷5="
Replace-Text "73÷2=" "63÷4="
Replace-Text "44÷4=" "83÷8="
Replace-Text "35÷7=" "26÷2="
Replace-Text "64÷7=" "86÷3="
Replace-Text "36÷3=" "97÷2="
Replace-Text "57÷7=" "37÷9="
Replace-Text "90÷4=" "74÷9="
Replace-Text "21÷6=" "51÷6="
Replace-Text "57÷6=" "55÷6="
Replace-Text "91÷3=" "42÷9="
Replace-Text "54÷5=" "33÷3="
Replace-Text "31÷2=" "32÷2="
Replace-Text "83÷6=" "22÷8="
Replace-Text "59÷8=" "25÷7="
Replace-Text "72÷3=" "80÷8="
Replace-Text "21÷4=" "59÷5="
Replace-Text "23÷8=" "39÷3="
Replace-Text "92÷2=" "52÷9="
Replace-Text "84÷9=" "72÷3="
Replace-Text "91÷9=" "12÷3="
Replace-Text "41÷2=" "37÷8="
Replace-Text "19÷4=" "85÷2="
Replace-Text "78÷4=" "77÷2="
Replace-Text "47÷4=" "61÷6="
